$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = 'Örebro University'
$ws.Range("B27").Value = 16
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 6.2
$ws.Range("E27").Value = 0.3
$ws.Range("F27").Value = 28.3

$ws.Range("A28").Value = 'Örebro University Hospital'
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 0
$ws.Range("D28:F28").ClearContents()

$ws.Range("A29").Value = 'Oslo University Hospital'
$ws.Range("B29").Value = 87
$ws.Range("C29").Value = 49
$ws.Range("D29").Value = 56.3
$ws.Range("E29").Value = 45.9
$ws.Range("F29").Value = 66.3

$ws.Range("A30").Value = 'Oulu University Hospital'
$ws.Range("B30").Value = 5
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 20
$ws.Range("E30").Value = 1
$ws.Range("F30").Value = 62.4

$ws.Range("A31").Value = 'Sahlgrenska University Hospital'
$ws.Range("B31").Value = 34
$ws.Range("C31").Value = 13
$ws.Range("D31").Value = 38.2
$ws.Range("E31").Value = 23.9
$ws.Range("F31").Value = 55.00000000000001

$ws.Range("A32").Value = 'Skane University Hospital'
$ws.Range("B32").Value = 14
$ws.Range("C32").Value = 6
$ws.Range("D32").Value = 42.9
$ws.Range("E32").Value = 21.4
$ws.Range("F32").Value = 67.40000000000001

$ws.Range("A33").Value = 'St. Olav’s University Hospital'
$ws.Range("B33").Value = 23
$ws.Range("C33").Value = 16
$ws.Range("D33").Value = 69.59999999999999
$ws.Range("E33").Value = 49.1
$ws.Range("F33").Value = 84.39999999999999

$ws.Range("A34").Value = 'Steno Diabetes Center Copenhagen'
$ws.Range("B34").Value = 7
$ws.Range("C34").Value = 4
$ws.Range("D34").Value = 57.1
$ws.Range("E34").Value = 25
$ws.Range("F34").Value = 84.2

$ws.Range("A35").Value = 'Stockholm South General Hospital'
$ws.Range("B35").Value = 3
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = 33.3
$ws.Range("E35").Value = 1.7
$ws.Range("F35").Value = 88.2

$ws.Range("A36").Value = 'Tampere University Hospital'
$ws.Range("B36").Value = 15
$ws.Range("C36").Value = 4
$ws.Range("D36").Value = 26.7
$ws.Range("E36").Value = 10.9
$ws.Range("F36").Value = 52

$ws.Range("A37").Value = 'The National University Hospital of Iceland'
$ws.Range("B37").Value = 2
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 50
$ws.Range("E37").Value = 2.6
$ws.Range("F37").Value = 97.39999999999999

$ws.Range("A38").Value = 'Turku University Hospital'
$ws.Range("B38").Value = 39
$ws.Range("C38").Value = 26
$ws.Range("D38").Value = 66.7
$ws.Range("E38").Value = 51
$ws.Range("F38").Value = 79.40000000000001

$ws.Range("A39").Value = 'UiT The Arctic University of Norway'
$ws.Range("B39").Value = 13
$ws.Range("C39").Value = 6
$ws.Range("D39").Value = 46.2
$ws.Range("E39").Value = 23.2
$ws.Range("F39").Value = 70.89999999999999

$ws.Range("A40").Value = 'Umeå University'
$ws.Range("B40").Value = 40
$ws.Range("C40").Value = 22
$ws.Range("D40").Value = 55
$ws.Range("E40").Value = 39.8
$ws.Range("F40").Value = 69.3

$ws.Range("A41").Value = 'University Hospital of North Norway'
$ws.Range("B41").Value = 17
$ws.Range("C41").Value = 12
$ws.Range("D41").Value = 70.59999999999999
$ws.Range("E41").Value = 46.9
$ws.Range("F41").Value = 86.7

$ws.Range("A42").Value = 'University Hospital of Umeå'
$ws.Range("B42").Value = 1
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 94.89999999999999

$ws.Range("A43").Value = 'University of Bergen'
$ws.Range("B43").Value = 30
$ws.Range("C43").Value = 13
$ws.Range("D43").Value = 43.3
$ws.Range("E43").Value = 27.4
$ws.Range("F43").Value = 60.8

$ws.Range("A44").Value = 'University of Copenhagen'
$ws.Range("B44").Value = 90
$ws.Range("C44").Value = 45
$ws.Range("D44").Value = 50
$ws.Range("E44").Value = 39.90000000000001
$ws.Range("F44").Value = 60.09999999999999

$ws.Range("A45").Value = 'University of Eastern Finland'
$ws.Range("B45").Value = 12
$ws.Range("C45").Value = 7
$ws.Range("D45").Value = 58.3
$ws.Range("E45").Value = 32
$ws.Range("F45").Value = 80.7

$ws.Range("A46").Value = 'University of Helsinki'
$ws.Range("B46").Value = 17
$ws.Range("C46").Value = 9
$ws.Range("D46").Value = 52.9
$ws.Range("E46").Value = 31
$ws.Range("F46").Value = 73.8

$ws.Range("A47").Value = 'University of Iceland'
$ws.Range("B47").Value = 5
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = 20
$ws.Range("E47").Value = 1
$ws.Range("F47").Value = 62.4

$ws.Range("A48").Value = 'University of Oslo'
$ws.Range("B48").Value = 21
$ws.Range("C48").Value = 12
$ws.Range("D48").Value = 57.1
$ws.Range("E48").Value = 36.5
$ws.Range("F48").Value = 75.5

$ws.Range("A49").Value = 'University of Oulu'
$ws.Range("B49").Value = 19
$ws.Range("C49").Value = 5
$ws.Range("D49").Value = 26.3
$ws.Range("E49").Value = 11.8
$ws.Range("F49").Value = 48.8

$ws.Range("A50").Value = 'University of Southern Denmark'
$ws.Range("B50").Value = 42
$ws.Range("C50").Value = 23
$ws.Range("D50").Value = 54.8
$ws.Range("E50").Value = 39.90000000000001
$ws.Range("F50").Value = 68.8

$ws.Range("A51").Value = 'University of Tampere'
$ws.Range("B51").Value = 8
$ws.Range("C51").Value = 1
$ws.Range("D51").Value = 12.5
$ws.Range("E51").Value = 0.6
$ws.Range("F51").Value = 47.09999999999999

$ws.Range("A52").Value = 'University of Turku'
$ws.Range("B52").Value = 16
$ws.Range("C52").Value = 7
$ws.Range("D52").Value = 43.8
$ws.Range("E52").Value = 23.1
$ws.Range("F52").Value = 66.8

$ws.Range("A53").Value = 'Uppsala Academic Hospital'
$ws.Range("B53").Value = 6
$ws.Range("C53").Value = 1
$ws.Range("D53").Value = 16.7
$ws.Range("E53").Value = 0.8999999999999999
$ws.Range("F53").Value = 56.39999999999999

$ws.Range("A54").Value = 'Uppsala University'
$ws.Range("B54").Value = 41
$ws.Range("C54").Value = 18
$ws.Range("D54").Value = 43.9
$ws.Range("E54").Value = 29.9
$ws.Range("F54").Value = 59

$ws.Range("A55").Value = 'Zealand University Hospital'
$ws.Range("B55").Value = 14
$ws.Range("C55").Value = 8
$ws.Range("D55").Value = 57.1
$ws.Range("E55").Value = 32.6
$ws.Range("F55").Value = 78.60000000000001
